$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Before:
#   Sheet 1: "总计"    - one summary row for 2022-Q3
#   Sheet 2: "2022-Q3" - per-fund detail for 2022-Q3 (3 funds)
#
# After:
#   Sheet 1: "总计"    - summary rows for 2022-Q4 (new, on top) then 2022-Q3
#   Sheet 2: "2022-Q4" - per-fund detail for 2022-Q4 (1 fund)  [was sheet 2]
#   Sheet 3: "2022-Q3" - per-fund detail for 2022-Q3 (3 funds) [duplicate of
#                         the original sheet 2, kept as an archive]
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# --- 1. Duplicate the 2022-Q3 detail sheet so the old figures are kept ----
# Placed immediately after the source sheet; inherits all data/styles/
# view-state (e.g. the active-tab flag) from it automatically.
$q3.Copy($null, $q3)
$q3Archive = $wb.Worksheets.Item(3)

# Rename: the original sheet becomes the (to-be-overwritten) 2022-Q4 sheet,
# the fresh duplicate keeps the 2022-Q3 name/data untouched.
$q3.Name = "2022-Q4"
$q3Archive.Name = "2022-Q3"

# --- 2. Rewrite the original sheet's data with the 2022-Q4 figures --------
# Only one fund this quarter, so rows 3:4 are no longer needed.
$q3.Range("A3:H4").EntireRow.Delete()

# Columns B-G hold text values even though several look numeric (fund code,
# size, position %, etc.) - force text formatting before writing so Excel
# doesn't silently coerce them to numbers, then drop the format again so no
# stray style index is left on the cells (matches the archive sheet).
$q4TextCells = $q3.Range("B2:G2")
$q4TextCells.NumberFormat = "@"

$q3.Range("B2").Value = "519029"
$q3.Range("C2").Value = "华夏稳增混合"
$q3.Range("D2").Value = "9.01"
$q3.Range("E2").Value = "94.55"
$q3.Range("F2").Value = "5.67"
$q3.Range("G2").Value = "0.5109"
$q3.Range("H2").Value = 4

$q4TextCells.ClearFormats()

# --- 3. Update the "总计" summary sheet ------------------------------------
# Existing row 2 (2022-Q3 totals) shifts down to row 3 unchanged; row 2 is
# overwritten with the new 2022-Q4 totals. Copy A2's formatting onto A3
# first so the shifted row keeps the same cell style as before.
$summary.Range("A2").Copy($summary.Range("A3"))

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.49

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.51
